# Adding more users + removing unused variable
#
# Row 20 (previously the first of a run of "Anonymous" placeholder rows)
# becomes a new named user "Serge".
# A brand-new row 25 is appended, duplicating the formatting of row 24,
# and is populated with "Anonymous" (keeping the same number of
# "Anonymous" placeholder rows as before, just shifted down by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 was "Anonymous" -- rename it to the new user "Serge".
$ws.Range("A20").Value = "Serge"

# Duplicate row 24 (values + formatting + row height) into the new row 25.
$ws.Range("A24:E24").Copy()
$ws.Range("A25:E25").PasteSpecial(-4122)
$ws.Rows.Item(25).RowHeight = $ws.Rows.Item(24).RowHeight()

# The new row keeps the "Anonymous" placeholder value.
$ws.Range("A25").Value = "Anonymous"
